$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3625542819499969
$ws.Range("B1").Value = 0.6507856845855713
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.37424099445343
$ws.Range("E1").Value = 0.8454117774963379
